# Replace the "<Project name>" placeholder in the Subtitle paragraph with
# "<kiwifeeds>", reproducing the exact run/proofErr structure Word leaves
# behind when the typed word is flagged by the spell checker: the text is
# split into three runs ("<", "kiwifeeds", ">") with spellStart/spellEnd
# proofErr markers bracketing the misspelled word.

$d = $word.ActiveDocument
$sel = $word.Selection

$found = $sel.Find.Execute("<Project name>", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $found) {
    throw "Could not find the '<Project name>' placeholder"
}

$rng = $sel.Range

# Clear the placeholder text first so InsertXML below replaces it in place
# rather than appending after it.
$rng.Text = ""

$xml = @'
<?xml version="1.0" encoding="UTF-8" standalone="yes"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">
  <pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml" pkg:padding="512">
    <pkg:xmlData>
      <w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml">
        <w:body>
          <w:p w14:paraId="5B5D1750" w14:textId="234AB331" w:rsidR="007E00E0" w:rsidRPr="00110419" w:rsidRDefault="007E00E0" w:rsidP="007E00E0">
            <w:pPr>
              <w:pStyle w:val="Subtitle"/>
              <w:rPr>
                <w:rFonts w:cs="Segoe UI"/>
              </w:rPr>
            </w:pPr>
            <w:r w:rsidRPr="00110419">
              <w:rPr>
                <w:rFonts w:cs="Segoe UI"/>
              </w:rPr>
              <w:t>&lt;</w:t>
            </w:r>
            <w:proofErr w:type="spellStart"/>
            <w:r>
              <w:rPr>
                <w:rFonts w:cs="Segoe UI"/>
              </w:rPr>
              <w:t>kiwifeeds</w:t>
            </w:r>
            <w:proofErr w:type="spellEnd"/>
            <w:r>
              <w:rPr>
                <w:rFonts w:cs="Segoe UI"/>
              </w:rPr>
              <w:t>&gt;</w:t>
            </w:r>
          </w:p>
        </w:body>
      </w:document>
    </pkg:xmlData>
  </pkg:part>
</pkg:package>
'@

$rng.InsertXML($xml) | Out-Null
